$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Replace the "TRUE" shared-string placeholder values in B2:E24 with real
# boolean FALSE values (displayed as FALSE, stored as t="b" / <v>0</v>).
$rng = $ws.Range("B2:E24")
$rng.Value = $false

# Update the active selection shown in the sheet view.
$ws.Activate()
$ws.Range("J18").Select()
